{"js": "// Soar-RL Tutorial: update the default exploration policy description.\n//\n// 1) The paragraph beginning \"Soar-RL's default exploration policy is\n//    epsilon-greedy, with an epsilon value of 0.1. ...\" is rewritten to\n//    explain that the new default policy is softmax, and that Soar-RL\n//    automatically switches to epsilon-greedy when first enabled.\n// 2) The paragraph \"Acceptable values for epsilon are numbers between 0\n//    and 1 ...\" gets \"epsilon\" italicized.\n// 3) A new empty paragraph is appended at the very end of the document.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the two paragraphs we need to rewrite by their (still unique)\n// original text, rather than hard-coded indices.\nlet policyParaIndex = -1;\nlet acceptableParaIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (policyParaIndex === -1 && t.indexOf(\"Soar-RL\\u2019s default exploration policy is\") === 0) {\n    policyParaIndex = i;\n  }\n  if (acceptableParaIndex === -1 && t.indexOf(\"Acceptable values for epsilon are numbers between 0 and 1\") === 0) {\n    acceptableParaIndex = i;\n  }\n}\nif (policyParaIndex === -1) {\n  throw new Error(\"Could not find the 'Soar-RL\\u2019s default exploration policy' paragraph.\");\n}\nif (acceptableParaIndex === -1) {\n  throw new Error(\"Could not find the 'Acceptable values for epsilon' paragraph.\");\n}\n\n// Helper: clear a paragraph, then rebuild it from a list of\n// [text, italic] pairs, always appending to the paragraph's end so each\n// new run starts from the paragraph's (non-italic) base formatting\n// instead of inheriting the italic formatting of the previous run.\nasync function rebuildParagraph(paragraph, parts) {\n  paragraph.clear();\n  await context.sync();\n  for (const [text, italic] of parts) {\n    const run = paragraph.insertText(text, Word.InsertLocation.end);\n    if (italic) {\n      run.font.italic = true;\n    }\n    await context.sync();\n  }\n}\n\nconst policyParagraph = paragraphs.items[policyParaIndex];\nawait rebuildParagraph(policyParagraph, [\n  [\"When Soar is first started, the default exploration policy is \", false],\n  [\"softmax\", true],\n  [\".  However, the first time Soar-RL is enabled, the architecture automatically changes the exploration policy to \", false],\n  [\"epsilon-greedy\", true],\n  [\", \", false],\n  [\"a policy more suitable for RL agents\", false],\n  [\".  \", false],\n  [\"The default value of \", false],\n  [\"epsilon\", true],\n  [\" is 0.1, dictating that \", false],\n  [\"90% of the time the operator with greatest numerical preference value is chosen, while the remaining 10% of the time a random selection is made from all acceptable proposed operators.  You can change the \", false],\n  [\"epsilon\", true],\n  [\" value by issuing the following command:\", false],\n]);\n\nconst acceptableParagraph = paragraphs.items[acceptableParaIndex];\nawait rebuildParagraph(acceptableParagraph, [\n  [\"Acceptable values for \", false],\n  [\"epsilon\", true],\n  [\" are numbers between 0 and 1\", false],\n  [\" (inclusive)\", false],\n  [\".  You may note, by the definition, that a value of 0 will eliminate the chance of exploration and a value \", false],\n  [\"of 1 will result in a uniformly\", false],\n  [\" random selection.\", false],\n]);\n\n// Append a new, empty paragraph at the very end of the document body.\nbody.insertParagraph(\"\", Word.InsertLocation.end);\nawait context.sync();\n", "ps1": "# Soar-RL Tutorial: update the default exploration policy description.\n#\n# 1) The paragraph beginning \"Soar-RL's default exploration policy is\n#    epsilon-greedy, with an epsilon value of 0.1. ...\" is rewritten to\n#    explain that the new default policy is softmax, and that Soar-RL\n#    automatically switches to epsilon-greedy when first enabled.\n# 2) The paragraph \"Acceptable values for epsilon are numbers between 0\n#    and 1 ...\" gets \"epsilon\" italicized.\n# 3) A new empty paragraph is appended at the very end of the document.\n\n$d = $word.ActiveDocument\n\n# Inserts $text at $position (a document character offset), optionally\n# italic, returning the offset just after the inserted text. Using a\n# freshly-collapsed Range for every insertion (rather than re-using /\n# extending a previously formatted Range) keeps each run's formatting\n# independent, so non-italic runs stay free of any explicit <w:i w:val=\"0\"/>.\nfunction Insert-Run($position, $text, $italic) {\n    $r = $d.Range($position, $position)\n    $r.InsertAfter($text)\n    $r.End = $r.Start + $text.Length\n    if ($italic) {\n        $r.Font.Italic = 1\n    }\n    return $r.End\n}\n\n# Clears the text of paragraph $paragraph (leaving its trailing paragraph\n# mark / paragraph-level formatting untouched) and returns the offset at\n# which new runs should be inserted.\nfunction Clear-ParagraphText($paragraph) {\n    $rng = $paragraph.Range\n    $start = $rng.Start\n    $d.Range($start, $rng.End - 1).Text = \"\"\n    return $start\n}\n\n# Locate the two paragraphs we need to rewrite by their (still unique)\n# original text, rather than a hard-coded paragraph index.\n$paragraphs = $d.Paragraphs\n$policyParagraph = $null\n$acceptableParagraph = $null\nfor ($i = 1; $i -le $paragraphs.Count; $i++) {\n    $para = $paragraphs.Item($i)\n    $t = $para.Range.Text\n    if (($null -eq $policyParagraph) -and $t.StartsWith(\"Soar-RL\u2019s default exploration policy is\")) {\n        $policyParagraph = $para\n    }\n    if (($null -eq $acceptableParagraph) -and $t.StartsWith(\"Acceptable values for epsilon are numbers between 0 and 1\")) {\n        $acceptableParagraph = $para\n    }\n}\nif ($null -eq $policyParagraph) {\n    throw \"Could not find the 'Soar-RL's default exploration policy' paragraph.\"\n}\nif ($null -eq $acceptableParagraph) {\n    throw \"Could not find the 'Acceptable values for epsilon' paragraph.\"\n}\n\n# --- Paragraph 1: default exploration policy explanation ---\n$pos = Clear-ParagraphText $policyParagraph\n$pos = Insert-Run $pos \"When Soar is first started, the default exploration policy is \" $false\n$pos = Insert-Run $pos \"softmax\" $true\n$pos = Insert-Run $pos \".  However, the first time Soar-RL is enabled, the architecture automatically changes the exploration policy to \" $false\n$pos = Insert-Run $pos \"epsilon-greedy\" $true\n$pos = Insert-Run $pos \", \" $false\n$pos = Insert-Run $pos \"a policy more suitable for RL agents\" $false\n$pos = Insert-Run $pos \".  \" $false\n$pos = Insert-Run $pos \"The default value of \" $false\n$pos = Insert-Run $pos \"epsilon\" $true\n$pos = Insert-Run $pos \" is 0.1, dictating that \" $false\n$pos = Insert-Run $pos \"90% of the time the operator with greatest numerical preference value is chosen, while the remaining 10% of the time a random selection is made from all acceptable proposed operators.  You can change the \" $false\n$pos = Insert-Run $pos \"epsilon\" $true\n$pos = Insert-Run $pos \" value by issuing the following command:\" $false\n\n# --- Paragraph 2: acceptable epsilon values ---\n$pos = Clear-ParagraphText $acceptableParagraph\n$pos = Insert-Run $pos \"Acceptable values for \" $false\n$pos = Insert-Run $pos \"epsilon\" $true\n$pos = Insert-Run $pos \" are numbers between 0 and 1\" $false\n$pos = Insert-Run $pos \" (inclusive)\" $false\n$pos = Insert-Run $pos \".  You may note, by the definition, that a value of 0 will eliminate the chance of exploration and a value \" $false\n$pos = Insert-Run $pos \"of 1 will result in a uniformly\" $false\n$pos = Insert-Run $pos \" random selection.\" $false\n\n# --- Append a new, empty paragraph at the very end of the document ---\n$endRange = $d.Content\n$endRange.Collapse(0) # wdCollapseEnd\n$endRange.InsertParagraphAfter()\n"}
